# "Refactorizando codigo para op"
# - Drop the workbook structure-protection record (no password was set).
# - On the "ruta_actual" sheet: write today's route date into B1 (next to
#   the "Ruta actual:" label) and clear out the stale date that used to
#   live down in A3, plus the leftover empty cell in I1.
# - On the "clientes" sheet: clear the empty placeholder left in G2.

$wb = $excel.ActiveWorkbook

$wb.Unprotect()

$wsRuta = $wb.Worksheets.Item("ruta_actual")
$wsRuta.Range("B1").Value = "20240706"
$wsRuta.Range("A3").ClearContents()
$wsRuta.Range("I1").ClearContents()
$wsRuta.Range("B1").Select() | Out-Null

$wsClientes = $wb.Worksheets.Item("clientes")
$wsClientes.Range("G2").ClearContents()
